$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.549.00"
$ws.Range("E2").Value = "  +3.85%  "

$ws.Range("D3").Value = "3.487.10"
$ws.Range("E3").Value = "  +2.38%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'591.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.69%  "

$ws.Range("D6").Value = "'169.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.02%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "3.484.03"
$ws.Range("E8").Value = "  +2.27%  "

$ws.Range("E9").Value = "  +7.98%  "

$ws.Range("D10").Value = "'7.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.19%  "

$ws.Range("E11").Value = "  +6.89%  "

$ws.Range("E12").Value = "  +4.02%  "

$ws.Range("D13").Value = "4.087.92"
$ws.Range("E13").Value = "  +2.30%  "

$ws.Range("D14").Value = "'0.135"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").Value = "'28.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.88%  "

$ws.Range("E16").Value = "  +3.37%  "

$ws.Range("D17").Value = "66.587.19"
$ws.Range("E17").Value = "  +3.86%  "

$ws.Range("D18").Value = "3.472.05"
$ws.Range("E18").Value = "  +2.90%  "

$ws.Range("D19").Value = "'6.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.26%  "

$ws.Range("D20").Value = "'14.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.96%  "

$ws.Range("D21").Value = "'391.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.21%  "

$ws.Range("D22").Value = "'7.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "

$ws.Range("D23").Value = "'72.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.94%  "

$ws.Range("E24").Value = "  -0.28%  "

$ws.Range("E25").Value = "  +4.80%  "

$ws.Range("E26").Value = "  +6.25%  "

$ws.Range("D27").Value = "'10.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.54%  "

$ws.Range("E28").Value = "  +1.41%  "

$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").Value = "'6.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.84%  "

$ws.Range("E31").Value = "  +5.09%  "

$ws.Range("E33").Value = "  +3.66%  "

$ws.Range("D34").Value = "'7.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.66%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("E36").Value = "  +8.60%  "

$ws.Range("D37").Value = "'161.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.31%  "

$ws.Range("D38").Value = "'0.891"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.11%  "

$ws.Range("E39").Value = "  +6.36%  "

$ws.Range("D40").Value = "'6.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.37%  "

$ws.Range("D41").Value = "'0.0744"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.19%  "

$ws.Range("D42").Value = "'26.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.78%  "

$ws.Range("E43").Value = "  +6.66%  "

$ws.Range("D44").Value = "'26.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.68%  "

$ws.Range("D45").Value = "'43.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.05%  "

$ws.Range("D46").Value = "2.765.50"
$ws.Range("E46").Value = "  +1.52%  "

$ws.Range("D47").Value = "'0.0313"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.93%  "

$ws.Range("E48").Value = "  +4.17%  "

$ws.Range("D49").Value = "'345.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.79%  "

$ws.Range("E50").Value = "  +4.69%  "

$ws.Range("D51").Value = "'33.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.63%  "
